$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '27.160.68'
$ws.Cells.Item(2, 5).Value = '  -1.87%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.562.39'
$ws.Cells.Item(3, 5).Value = '  -1.63%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.10%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '206.34'
$ws.Cells.Item(5, 5).Value = '  -0.36%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -1.63%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.10%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '22.27'
$ws.Cells.Item(8, 5).Value = '  +0.22%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -1.96%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +0.04%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.90%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.783.86'
$ws.Cells.Item(12, 5).Value = '  -1.68%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.558.50'
$ws.Cells.Item(13, 5).Value = '  -1.81%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -2.14%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '0.515'
$ws.Cells.Item(15, 5).Value = '  -2.72%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '62.92'
$ws.Cells.Item(16, 5).Value = '  -0.76%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '27.126.75'
$ws.Cells.Item(17, 5).Value = '  -1.96%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '213.27'
$ws.Cells.Item(18, 5).Value = '  -2.80%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '0.0₃0688'
$ws.Cells.Item(19, 5).Value = '  -1.03%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '7.22'
$ws.Cells.Item(20, 5).Value = '  -1.29%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.05%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '4.11'
$ws.Cells.Item(22, 5).Value = '  -0.69%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '9.38'
$ws.Cells.Item(23, 5).Value = '  -2.21%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '1.97'
$ws.Cells.Item(24, 5).Value = '  +0.38%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '152.15'
$ws.Cells.Item(25, 5).Value = '  -0.83%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '6.58'
$ws.Cells.Item(26, 5).Value = '  -3.85%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '14.87'
$ws.Cells.Item(27, 5).Value = '  -1.66%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -0.07%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -1.60%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -1.00%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '0.0464'
$ws.Cells.Item(31, 5).Value = '  -0.74%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -1.76%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '1.380.46'
$ws.Cells.Item(33, 5).Value = '  +0.81%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +0.62%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '1.55'
$ws.Cells.Item(35, 5).Value = '  +0.68%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'HuobiToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(36, 4).Value = '2.28'
$ws.Cells.Item(36, 5).Value = '  -1.06%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.940'
$ws.Cells.Item(37, 5).Value = '  -4.11%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -1.34%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.814'
$ws.Cells.Item(39, 5).Value = '  -1.23%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '0.518'
$ws.Cells.Item(40, 5).Value = '  -3.43%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.03%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '0.993'
$ws.Cells.Item(42, 5).Value = '  +1.96%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '1.78'
$ws.Cells.Item(43, 5).Value = '  +2.80%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Aave'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(44, 4).Value = '63.44'
$ws.Cells.Item(44, 5).Value = '  -1.03%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'MXToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(45, 4).Value = '2.17'
$ws.Cells.Item(45, 5).Value = '  -0.02%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '5.22'
$ws.Cells.Item(46, 5).Value = '  -0.16%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '1.696.86'
$ws.Cells.Item(47, 5).Value = '  -1.66%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '85.55'
$ws.Cells.Item(48, 5).Value = '  -2.28%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '0.0₇0999'
$ws.Cells.Item(49, 5).Value = '  -0.79%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -0.74%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -0.05%  '
